$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "Inflammatory-Mac"
$ws.Range("G2").Value = 169.5951436666667
$ws.Range("H2").Value = 508.785431
$ws.Range("I2").Value = 0.2074259764082431
$ws.Range("J2").Value = 0.2074259764082431
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05743166666666667
$ws.Range("N2").Value = 0.172295
$ws.Range("O2").Value = 0.4025715794441874
$ws.Range("P2").Value = 0.4025715794441875
$ws.Range("Q2").Value = 9.740131759349445
$ws.Range("R2").Value = 87.661185834145
$ws.Range("S2").Value = 0.08350380294041916
$ws.Range("T2").Value = 0.08350380294041919

# Row 3
$ws.Range("G3").Value = 169.5951436666667
$ws.Range("H3").Value = 508.785431
$ws.Range("I3").Value = 0.2074259764082431
$ws.Range("J3").Value = 0.2074259764082431
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.08523033333333334
$ws.Range("N3").Value = 0.255691
$ws.Range("O3").Value = 0.5974284205558126
$ws.Range("P3").Value = 0.5974284205558126
$ws.Range("Q3").Value = 14.45465062642456
$ws.Range("R3").Value = 130.091855637821
$ws.Range("S3").Value = 0.1239221734678239
$ws.Range("T3").Value = 0.1239221734678239

# Row 4
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("I4").Value = 0.4056457917095931
$ws.Range("J4").Value = 0.405645791709593
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.05743166666666667
$ws.Range("N4").Value = 0.172295
$ws.Range("O4").Value = 0.4025715794441874
$ws.Range("P4").Value = 0.4025715794441875
$ws.Range("Q4").Value = 19.04796847189889
$ws.Range("R4").Value = 171.43171624709
$ws.Range("S4").Value = 0.1633014670634187
$ws.Range("T4").Value = 0.1633014670634187

# Row 5
$ws.Range("I5").Value = 0.4056457917095931
$ws.Range("J5").Value = 0.405645791709593
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.08523033333333334
$ws.Range("N5").Value = 0.255691
$ws.Range("O5").Value = 0.5974284205558126
$ws.Range("P5").Value = 0.5974284205558126
$ws.Range("Q5").Value = 28.26776230620911
$ws.Range("R5").Value = 254.409860755882
$ws.Range("S5").Value = 0.2423443246461743
$ws.Range("T5").Value = 0.2423443246461743

# Row 6
$ws.Range("D6").Value = "Inflammatory-Mac"
$ws.Range("G6").Value = 98.17454766666667
$ws.Range("H6").Value = 294.523643
$ws.Range("I6").Value = 0.1200739064098473
$ws.Range("J6").Value = 0.1200739064098473
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.05743166666666667
$ws.Range("N6").Value = 0.172295
$ws.Range("O6").Value = 0.4025715794441874
$ws.Range("P6").Value = 0.4025715794441875
$ws.Range("Q6").Value = 5.638327896742778
$ws.Range("R6").Value = 50.744951070685
$ws.Range("S6").Value = 0.04833834215344576
$ws.Range("T6").Value = 0.04833834215344576

# Row 7
$ws.Range("G7").Value = 98.17454766666667
$ws.Range("H7").Value = 294.523643
$ws.Range("I7").Value = 0.1200739064098473
$ws.Range("J7").Value = 0.1200739064098473
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.08523033333333334
$ws.Range("N7").Value = 0.255691
$ws.Range("O7").Value = 0.5974284205558126
$ws.Range("P7").Value = 0.5974284205558126
$ws.Range("Q7").Value = 8.367449422479224
$ws.Range("R7").Value = 75.307044802313
$ws.Range("S7").Value = 0.07173556425640153
$ws.Range("T7").Value = 0.07173556425640151

# Row 8
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("G8").Value = 105.935201
$ws.Range("H8").Value = 317.805603
$ws.Range("I8").Value = 0.1295656940897851
$ws.Range("J8").Value = 0.1295656940897851
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.05743166666666667
$ws.Range("N8").Value = 0.172295
$ws.Range("O8").Value = 0.4025715794441874
$ws.Range("P8").Value = 0.4025715794441875
$ws.Range("Q8").Value = 6.084035152098334
$ws.Range("R8").Value = 54.756316368885
$ws.Range("S8").Value = 0.05215946611150721
$ws.Range("T8").Value = 0.05215946611150721

# Row 9
$ws.Range("G9").Value = 105.935201
$ws.Range("H9").Value = 317.805603
$ws.Range("I9").Value = 0.1295656940897851
$ws.Range("J9").Value = 0.1295656940897851
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.08523033333333334
$ws.Range("N9").Value = 0.255691
$ws.Range("O9").Value = 0.5974284205558126
$ws.Range("P9").Value = 0.5974284205558126
$ws.Range("Q9").Value = 9.028892492963667
$ws.Range("R9").Value = 81.260032436673
$ws.Range("S9").Value = 0.07740622797827791
$ws.Range("T9").Value = 0.0774062279782779

# Row 10
$ws.Range("D10").Value = "Inflammatory-Mac"
$ws.Range("G10").Value = 112.249611
$ws.Range("H10").Value = 336.748833
$ws.Range("I10").Value = 0.1372886313825315
$ws.Range("J10").Value = 0.1372886313825314
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.05743166666666667
$ws.Range("N10").Value = 0.172295
$ws.Range("O10").Value = 0.4025715794441874
$ws.Range("P10").Value = 0.4025715794441875
$ws.Range("Q10").Value = 6.446682242415
$ws.Range("R10").Value = 58.020140181735
$ws.Range("S10").Value = 0.05526850117539652
$ws.Range("T10").Value = 0.05526850117539653

# Row 11
$ws.Range("G11").Value = 112.249611
$ws.Range("H11").Value = 336.748833
$ws.Range("I11").Value = 0.1372886313825315
$ws.Range("J11").Value = 0.1372886313825314
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.08523033333333334
$ws.Range("N11").Value = 0.255691
$ws.Range("O11").Value = 0.5974284205558126
$ws.Range("P11").Value = 0.5974284205558126
$ws.Range("Q11").Value = 9.567071762067
$ws.Range("R11").Value = 86.103645858603
$ws.Range("S11").Value = 0.08202013020713494
$ws.Range("T11").Value = 0.08202013020713493
